$wb = $excel.ActiveWorkbook

# Target column widths (as stored in the worksheet XML "width" attribute)
# are 29.9777047293527 and 13.7470528738839. Excel's ColumnWidth COM
# property is specified in characters and gets snapped to the nearest
# pixel when written back out, so the COM-side values below are chosen
# to land in the middle of the pixel bucket that produces the closest
# possible stored width to the targets above.
$wideColumnWidth = 29.16666666666665
$narrowColumnWidth = 12.83333333333335

# --- Sheet "Overview" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws1.Columns.Item(5).ColumnWidth = $wideColumnWidth
$ws1.Columns.Item(6).ColumnWidth = $wideColumnWidth

# --- Sheet "zh-cn" ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("K2").Value = "2016-08-31 00:52:51"
# Assigning a bare "" clears/removes the cell entirely. Using a leading
# quote (text prefix) keeps the cell present with an empty string value,
# then resetting the style keeps it on the default "Normal" style (no
# quote-prefix indicator left behind on the cell).
$ws2.Range("P2").Value = "'"
$ws2.Range("P2").Style = "Normal"
$ws2.Columns.Item(3).ColumnWidth = $wideColumnWidth
$ws2.Columns.Item(16).ColumnWidth = $narrowColumnWidth

# --- Sheet "de-de" ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("K2").Value = "2016-08-31 00:52:58"
$ws3.Range("P2").Value = "'"
$ws3.Range("P2").Style = "Normal"
$ws3.Columns.Item(3).ColumnWidth = $wideColumnWidth
$ws3.Columns.Item(16).ColumnWidth = $narrowColumnWidth
